# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# to reflect the refreshed data snapshot committed at 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1825
    4  = 109
    7  = 1498
    9  = 589
    11 = 95
    12 = 16
    13 = 87
    14 = 219
    17 = 102
    18 = 120
    19 = 3545
    20 = 423
    21 = 315
    22 = 482
    23 = 108
    24 = 336
    25 = 27
    26 = 1291
    27 = 135
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
